$wb = $excel.ActiveWorkbook

# ==== Sheet 'Вопросы': rename the three questions, update creation date ====
$wsQuestions = $wb.Worksheets.Item("Вопросы")
$wsQuestions.Cells.Item(2, 3).Value = 'Первый вопрос'
$wsQuestions.Cells.Item(2, 4).Value = 44165
$wsQuestions.Cells.Item(3, 3).Value = 'Второй вопрос'
$wsQuestions.Cells.Item(3, 4).Value = 44165
$wsQuestions.Cells.Item(4, 3).Value = 'Третий вопрос'
$wsQuestions.Cells.Item(4, 4).Value = 44165

# ==== Sheet 'Ответы': extend with new PVI filter rows, update existing rows ====
$wsAnswers = $wb.Worksheets.Item("Ответы")

# Copy formatting (index column border/bold style; date number format) down to new rows 14-37
for ($r = 14; $r -le 37; $r++) {
    $wsAnswers.Range("A2").Copy($wsAnswers.Cells.Item($r, 1))
    $wsAnswers.Range("F2").Copy($wsAnswers.Cells.Item($r, 6))
}

# Set cell values for rows 2-37
$wsAnswers.Cells.Item(2, 1).Value = 0
$wsAnswers.Cells.Item(2, 2).Value = 78035
$wsAnswers.Cells.Item(2, 3).Value = 'Первый вопрос'
$wsAnswers.Cells.Item(2, 4).Value = 1
$wsAnswers.Cells.Item(2, 5).Value = 'None'
$wsAnswers.Cells.Item(2, 6).Value = 44165
$wsAnswers.Cells.Item(2, 7).Value = 'Н-Уренгойское ЛПУМГ (ПВП № 1, КС Пуртазовская)'
$wsAnswers.Cells.Item(3, 1).Value = 1
$wsAnswers.Cells.Item(3, 2).Value = 78035
$wsAnswers.Cells.Item(3, 3).Value = 'Второй вопрос'
$wsAnswers.Cells.Item(3, 4).Value = 1
$wsAnswers.Cells.Item(3, 5).Value = 'None'
$wsAnswers.Cells.Item(3, 6).Value = 44165
$wsAnswers.Cells.Item(3, 7).Value = 'Н-Уренгойское ЛПУМГ (ПВП № 1, КС Пуртазовская)'
$wsAnswers.Cells.Item(4, 1).Value = 2
$wsAnswers.Cells.Item(4, 2).Value = 78035
$wsAnswers.Cells.Item(4, 3).Value = 'Третий вопрос'
$wsAnswers.Cells.Item(4, 4).Value = 1
$wsAnswers.Cells.Item(4, 5).Value = 'None'
$wsAnswers.Cells.Item(4, 6).Value = 44165
$wsAnswers.Cells.Item(4, 7).Value = 'Н-Уренгойское ЛПУМГ (ПВП № 1, КС Пуртазовская)'
$wsAnswers.Cells.Item(5, 1).Value = 3
$wsAnswers.Cells.Item(5, 2).Value = 142820
$wsAnswers.Cells.Item(5, 3).Value = 'Первый вопрос'
$wsAnswers.Cells.Item(5, 4).Value = 1
$wsAnswers.Cells.Item(5, 5).Value = 'None'
$wsAnswers.Cells.Item(5, 6).Value = 44165
$wsAnswers.Cells.Item(5, 7).Value = 'Н-Уренгойское ЛПУМГ (ПВП № 1, КС Пуртазовская)'
$wsAnswers.Cells.Item(6, 1).Value = 4
$wsAnswers.Cells.Item(6, 2).Value = 142820
$wsAnswers.Cells.Item(6, 3).Value = 'Второй вопрос'
$wsAnswers.Cells.Item(6, 4).Value = 1
$wsAnswers.Cells.Item(6, 5).Value = 'None'
$wsAnswers.Cells.Item(6, 6).Value = 44165
$wsAnswers.Cells.Item(6, 7).Value = 'Н-Уренгойское ЛПУМГ (ПВП № 1, КС Пуртазовская)'
$wsAnswers.Cells.Item(7, 1).Value = 5
$wsAnswers.Cells.Item(7, 2).Value = 142820
$wsAnswers.Cells.Item(7, 3).Value = 'Третий вопрос'
$wsAnswers.Cells.Item(7, 4).Value = 1
$wsAnswers.Cells.Item(7, 5).Value = 'None'
$wsAnswers.Cells.Item(7, 6).Value = 44165
$wsAnswers.Cells.Item(7, 7).Value = 'Н-Уренгойское ЛПУМГ (ПВП № 1, КС Пуртазовская)'
$wsAnswers.Cells.Item(8, 1).Value = 6
$wsAnswers.Cells.Item(8, 2).Value = 150831
$wsAnswers.Cells.Item(8, 3).Value = 'Первый вопрос'
$wsAnswers.Cells.Item(8, 4).Value = 1
$wsAnswers.Cells.Item(8, 5).Value = 'None'
$wsAnswers.Cells.Item(8, 6).Value = 44165
$wsAnswers.Cells.Item(8, 7).Value = 'Н-Уренгойское ЛПУМГ (ПВП №2)'
$wsAnswers.Cells.Item(9, 1).Value = 7
$wsAnswers.Cells.Item(9, 2).Value = 150831
$wsAnswers.Cells.Item(9, 3).Value = 'Второй вопрос'
$wsAnswers.Cells.Item(9, 4).Value = 1
$wsAnswers.Cells.Item(9, 5).Value = 'None'
$wsAnswers.Cells.Item(9, 6).Value = 44165
$wsAnswers.Cells.Item(9, 7).Value = 'Н-Уренгойское ЛПУМГ (ПВП №2)'
$wsAnswers.Cells.Item(10, 1).Value = 8
$wsAnswers.Cells.Item(10, 2).Value = 150831
$wsAnswers.Cells.Item(10, 3).Value = 'Третий вопрос'
$wsAnswers.Cells.Item(10, 4).Value = 1
$wsAnswers.Cells.Item(10, 5).Value = 'None'
$wsAnswers.Cells.Item(10, 6).Value = 44165
$wsAnswers.Cells.Item(10, 7).Value = 'Н-Уренгойское ЛПУМГ (ПВП №2)'
$wsAnswers.Cells.Item(11, 1).Value = 9
$wsAnswers.Cells.Item(11, 2).Value = 9022
$wsAnswers.Cells.Item(11, 3).Value = 'Первый вопрос'
$wsAnswers.Cells.Item(11, 4).Value = 1
$wsAnswers.Cells.Item(11, 5).Value = 'None'
$wsAnswers.Cells.Item(11, 6).Value = 44165
$wsAnswers.Cells.Item(11, 7).Value = 'Пурпейское ЛПУМГ (КС - 01, Общежитие на 100 мест Ягенетская п/п)'
$wsAnswers.Cells.Item(12, 1).Value = 10
$wsAnswers.Cells.Item(12, 2).Value = 9022
$wsAnswers.Cells.Item(12, 3).Value = 'Второй вопрос'
$wsAnswers.Cells.Item(12, 4).Value = 1
$wsAnswers.Cells.Item(12, 5).Value = 'None'
$wsAnswers.Cells.Item(12, 6).Value = 44165
$wsAnswers.Cells.Item(12, 7).Value = 'Пурпейское ЛПУМГ (КС - 01, Общежитие на 100 мест Ягенетская п/п)'
$wsAnswers.Cells.Item(13, 1).Value = 11
$wsAnswers.Cells.Item(13, 2).Value = 9022
$wsAnswers.Cells.Item(13, 3).Value = 'Третий вопрос'
$wsAnswers.Cells.Item(13, 4).Value = 1
$wsAnswers.Cells.Item(13, 5).Value = 'None'
$wsAnswers.Cells.Item(13, 6).Value = 44165
$wsAnswers.Cells.Item(13, 7).Value = 'Пурпейское ЛПУМГ (КС - 01, Общежитие на 100 мест Ягенетская п/п)'
$wsAnswers.Cells.Item(14, 1).Value = 12
$wsAnswers.Cells.Item(14, 2).Value = 218044
$wsAnswers.Cells.Item(14, 3).Value = 'Первый вопрос'
$wsAnswers.Cells.Item(14, 4).Value = 1
$wsAnswers.Cells.Item(14, 5).Value = 'None'
$wsAnswers.Cells.Item(14, 6).Value = 44165
$wsAnswers.Cells.Item(14, 7).Value = 'Демьянское ЛПУМГ (КС - 7, Общежитие в п. Демьянка)'
$wsAnswers.Cells.Item(15, 1).Value = 13
$wsAnswers.Cells.Item(15, 2).Value = 218044
$wsAnswers.Cells.Item(15, 3).Value = 'Второй вопрос'
$wsAnswers.Cells.Item(15, 4).Value = 1
$wsAnswers.Cells.Item(15, 5).Value = 'None'
$wsAnswers.Cells.Item(15, 6).Value = 44165
$wsAnswers.Cells.Item(15, 7).Value = 'Демьянское ЛПУМГ (КС - 7, Общежитие в п. Демьянка)'
$wsAnswers.Cells.Item(16, 1).Value = 14
$wsAnswers.Cells.Item(16, 2).Value = 218044
$wsAnswers.Cells.Item(16, 3).Value = 'Третий вопрос'
$wsAnswers.Cells.Item(16, 4).Value = 1
$wsAnswers.Cells.Item(16, 5).Value = 'None'
$wsAnswers.Cells.Item(16, 6).Value = 44165
$wsAnswers.Cells.Item(16, 7).Value = 'Демьянское ЛПУМГ (КС - 7, Общежитие в п. Демьянка)'
$wsAnswers.Cells.Item(17, 1).Value = 15
$wsAnswers.Cells.Item(17, 2).Value = 233672
$wsAnswers.Cells.Item(17, 3).Value = 'Первый вопрос'
$wsAnswers.Cells.Item(17, 4).Value = 1
$wsAnswers.Cells.Item(17, 5).Value = 'None'
$wsAnswers.Cells.Item(17, 6).Value = 44165
$wsAnswers.Cells.Item(17, 7).Value = 'Сургутское ЛПУМГ (УСС "Факел")'
$wsAnswers.Cells.Item(18, 1).Value = 16
$wsAnswers.Cells.Item(18, 2).Value = 233672
$wsAnswers.Cells.Item(18, 3).Value = 'Второй вопрос'
$wsAnswers.Cells.Item(18, 4).Value = 1
$wsAnswers.Cells.Item(18, 5).Value = 'None'
$wsAnswers.Cells.Item(18, 6).Value = 44165
$wsAnswers.Cells.Item(18, 7).Value = 'Сургутское ЛПУМГ (УСС "Факел")'
$wsAnswers.Cells.Item(19, 1).Value = 17
$wsAnswers.Cells.Item(19, 2).Value = 233672
$wsAnswers.Cells.Item(19, 3).Value = 'Третий вопрос'
$wsAnswers.Cells.Item(19, 4).Value = 1
$wsAnswers.Cells.Item(19, 5).Value = 'None'
$wsAnswers.Cells.Item(19, 6).Value = 44165
$wsAnswers.Cells.Item(19, 7).Value = 'Сургутское ЛПУМГ (УСС "Факел")'
$wsAnswers.Cells.Item(20, 1).Value = 18
$wsAnswers.Cells.Item(20, 2).Value = 265177
$wsAnswers.Cells.Item(20, 3).Value = 'Первый вопрос'
$wsAnswers.Cells.Item(20, 4).Value = 5
$wsAnswers.Cells.Item(20, 5).Value = 'None'
$wsAnswers.Cells.Item(20, 6).Value = 44165
$wsAnswers.Cells.Item(20, 7).Value = 'Сургутское ЛПУМГ (УСС "Факел")'
$wsAnswers.Cells.Item(21, 1).Value = 19
$wsAnswers.Cells.Item(21, 2).Value = 265177
$wsAnswers.Cells.Item(21, 3).Value = 'Второй вопрос'
$wsAnswers.Cells.Item(21, 4).Value = 4
$wsAnswers.Cells.Item(21, 5).Value = 'None'
$wsAnswers.Cells.Item(21, 6).Value = 44165
$wsAnswers.Cells.Item(21, 7).Value = 'Сургутское ЛПУМГ (УСС "Факел")'
$wsAnswers.Cells.Item(22, 1).Value = 20
$wsAnswers.Cells.Item(22, 2).Value = 265177
$wsAnswers.Cells.Item(22, 3).Value = 'Третий вопрос'
$wsAnswers.Cells.Item(22, 4).Value = 3
$wsAnswers.Cells.Item(22, 5).Value = 'None'
$wsAnswers.Cells.Item(22, 6).Value = 44165
$wsAnswers.Cells.Item(22, 7).Value = 'Сургутское ЛПУМГ (УСС "Факел")'
$wsAnswers.Cells.Item(23, 1).Value = 21
$wsAnswers.Cells.Item(23, 2).Value = 292964
$wsAnswers.Cells.Item(23, 3).Value = 'Первый вопрос'
$wsAnswers.Cells.Item(23, 4).Value = 1
$wsAnswers.Cells.Item(23, 5).Value = 'None'
$wsAnswers.Cells.Item(23, 6).Value = 44167
$wsAnswers.Cells.Item(23, 7).Value = 'Н-Уренгойское ЛПУМГ (ПВП № 1, КС Пуртазовская)'
$wsAnswers.Cells.Item(24, 1).Value = 22
$wsAnswers.Cells.Item(24, 2).Value = 292964
$wsAnswers.Cells.Item(24, 3).Value = 'Второй вопрос'
$wsAnswers.Cells.Item(24, 4).Value = 1
$wsAnswers.Cells.Item(24, 5).Value = 'None'
$wsAnswers.Cells.Item(24, 6).Value = 44167
$wsAnswers.Cells.Item(24, 7).Value = 'Н-Уренгойское ЛПУМГ (ПВП № 1, КС Пуртазовская)'
$wsAnswers.Cells.Item(25, 1).Value = 23
$wsAnswers.Cells.Item(25, 2).Value = 292964
$wsAnswers.Cells.Item(25, 3).Value = 'Третий вопрос'
$wsAnswers.Cells.Item(25, 4).Value = 1
$wsAnswers.Cells.Item(25, 5).Value = 'None'
$wsAnswers.Cells.Item(25, 6).Value = 44167
$wsAnswers.Cells.Item(25, 7).Value = 'Н-Уренгойское ЛПУМГ (ПВП № 1, КС Пуртазовская)'
$wsAnswers.Cells.Item(26, 1).Value = 24
$wsAnswers.Cells.Item(26, 2).Value = 137897
$wsAnswers.Cells.Item(26, 3).Value = 'Первый вопрос'
$wsAnswers.Cells.Item(26, 4).Value = 3
$wsAnswers.Cells.Item(26, 5).Value = 'None'
$wsAnswers.Cells.Item(26, 6).Value = 44167
$wsAnswers.Cells.Item(26, 7).Value = 'Н-Уренгойское ЛПУМГ (ПВП № 1, КС Пуртазовская)'
$wsAnswers.Cells.Item(27, 1).Value = 25
$wsAnswers.Cells.Item(27, 2).Value = 137897
$wsAnswers.Cells.Item(27, 3).Value = 'Второй вопрос'
$wsAnswers.Cells.Item(27, 4).Value = 2
$wsAnswers.Cells.Item(27, 5).Value = 'None'
$wsAnswers.Cells.Item(27, 6).Value = 44167
$wsAnswers.Cells.Item(27, 7).Value = 'Н-Уренгойское ЛПУМГ (ПВП № 1, КС Пуртазовская)'
$wsAnswers.Cells.Item(28, 1).Value = 26
$wsAnswers.Cells.Item(28, 2).Value = 137897
$wsAnswers.Cells.Item(28, 3).Value = 'Третий вопрос'
$wsAnswers.Cells.Item(28, 4).Value = 3
$wsAnswers.Cells.Item(28, 5).Value = 'None'
$wsAnswers.Cells.Item(28, 6).Value = 44167
$wsAnswers.Cells.Item(28, 7).Value = 'Н-Уренгойское ЛПУМГ (ПВП № 1, КС Пуртазовская)'
$wsAnswers.Cells.Item(29, 1).Value = 27
$wsAnswers.Cells.Item(29, 2).Value = 177809
$wsAnswers.Cells.Item(29, 3).Value = 'Первый вопрос'
$wsAnswers.Cells.Item(29, 4).Value = 3
$wsAnswers.Cells.Item(29, 5).Value = 'None'
$wsAnswers.Cells.Item(29, 6).Value = 44167
$wsAnswers.Cells.Item(29, 7).Value = 'Н-Уренгойское ЛПУМГ (ПВП № 1, КС Пуртазовская)'
$wsAnswers.Cells.Item(30, 1).Value = 28
$wsAnswers.Cells.Item(30, 2).Value = 177809
$wsAnswers.Cells.Item(30, 3).Value = 'Второй вопрос'
$wsAnswers.Cells.Item(30, 4).Value = 3
$wsAnswers.Cells.Item(30, 5).Value = 'None'
$wsAnswers.Cells.Item(30, 6).Value = 44167
$wsAnswers.Cells.Item(30, 7).Value = 'Н-Уренгойское ЛПУМГ (ПВП № 1, КС Пуртазовская)'
$wsAnswers.Cells.Item(31, 1).Value = 29
$wsAnswers.Cells.Item(31, 2).Value = 177809
$wsAnswers.Cells.Item(31, 3).Value = 'Третий вопрос'
$wsAnswers.Cells.Item(31, 4).Value = 3
$wsAnswers.Cells.Item(31, 5).Value = 'None'
$wsAnswers.Cells.Item(31, 6).Value = 44167
$wsAnswers.Cells.Item(31, 7).Value = 'Н-Уренгойское ЛПУМГ (ПВП № 1, КС Пуртазовская)'
$wsAnswers.Cells.Item(32, 1).Value = 30
$wsAnswers.Cells.Item(32, 2).Value = 238761
$wsAnswers.Cells.Item(32, 3).Value = 'Первый вопрос'
$wsAnswers.Cells.Item(32, 4).Value = 5
$wsAnswers.Cells.Item(32, 5).Value = 'None'
$wsAnswers.Cells.Item(32, 6).Value = 44167
$wsAnswers.Cells.Item(32, 7).Value = 'Н-Уренгойское ЛПУМГ (ПВП № 1, КС Пуртазовская)'
$wsAnswers.Cells.Item(33, 1).Value = 31
$wsAnswers.Cells.Item(33, 2).Value = 238761
$wsAnswers.Cells.Item(33, 3).Value = 'Второй вопрос'
$wsAnswers.Cells.Item(33, 4).Value = 5
$wsAnswers.Cells.Item(33, 5).Value = 'None'
$wsAnswers.Cells.Item(33, 6).Value = 44167
$wsAnswers.Cells.Item(33, 7).Value = 'Н-Уренгойское ЛПУМГ (ПВП № 1, КС Пуртазовская)'
$wsAnswers.Cells.Item(34, 1).Value = 32
$wsAnswers.Cells.Item(34, 2).Value = 238761
$wsAnswers.Cells.Item(34, 3).Value = 'Третий вопрос'
$wsAnswers.Cells.Item(34, 4).Value = 3
$wsAnswers.Cells.Item(34, 5).Value = 'None'
$wsAnswers.Cells.Item(34, 6).Value = 44167
$wsAnswers.Cells.Item(34, 7).Value = 'Н-Уренгойское ЛПУМГ (ПВП № 1, КС Пуртазовская)'
$wsAnswers.Cells.Item(35, 1).Value = 33
$wsAnswers.Cells.Item(35, 2).Value = 198786
$wsAnswers.Cells.Item(35, 3).Value = 'Первый вопрос'
$wsAnswers.Cells.Item(35, 4).Value = 5
$wsAnswers.Cells.Item(35, 5).Value = 'None'
$wsAnswers.Cells.Item(35, 6).Value = 44167
$wsAnswers.Cells.Item(35, 7).Value = 'Н-Уренгойское ЛПУМГ (ПВП № 1, КС Пуртазовская)'
$wsAnswers.Cells.Item(36, 1).Value = 34
$wsAnswers.Cells.Item(36, 2).Value = 198786
$wsAnswers.Cells.Item(36, 3).Value = 'Второй вопрос'
$wsAnswers.Cells.Item(36, 4).Value = 5
$wsAnswers.Cells.Item(36, 5).Value = 'None'
$wsAnswers.Cells.Item(36, 6).Value = 44167
$wsAnswers.Cells.Item(36, 7).Value = 'Н-Уренгойское ЛПУМГ (ПВП № 1, КС Пуртазовская)'
$wsAnswers.Cells.Item(37, 1).Value = 35
$wsAnswers.Cells.Item(37, 2).Value = 198786
$wsAnswers.Cells.Item(37, 3).Value = 'Третий вопрос'
$wsAnswers.Cells.Item(37, 4).Value = 3
$wsAnswers.Cells.Item(37, 5).Value = 'None'
$wsAnswers.Cells.Item(37, 6).Value = 44167
$wsAnswers.Cells.Item(37, 7).Value = 'Н-Уренгойское ЛПУМГ (ПВП № 1, КС Пуртазовская)'

# ==== Sheet 'Комментарии': remove sample data rows, keep header only ====
$wsComments = $wb.Worksheets.Item("Комментарии")
$wsComments.Range("A2:C3").Delete()
